# Update the table style applied to the three data tables (slides 14-16)
# from the default "No Style, Table Grid" variant to the new style
# {1366D48C-6E8F-4151-BC14-ED94205402D1}.

$p = $ppt.ActivePresentation

$newStyleId = "{1366D48C-6E8F-4151-BC14-ED94205402D1}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
